# Update the "取得日時" (retrieved datetime) column for all data rows
# on the "ランサーズ" sheet from 2025-12-06 12:44:18 to 2025-12-06 18:23:36.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ランサーズ")

$newTimestamp = "2025-12-06 18:23:36"

for ($row = 2; $row -le 8; $row++) {
    $ws.Cells.Item($row, 1).Value = $newTimestamp
}
